# Add a "Save" column (H) to the s_vals sheet, mirroring the style of the
# existing header row and leaving the data cells with the default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy style from G1 ("sum" header) then set text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column, rows 2-18.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 0
    18 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
